# Apply cryptos.xlsx price/volume update
# (Updated cryptos list on Sun May 21 16:55:55 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.132.23"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.823.91"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'312.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("D8").Value = "'0.3625"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").Value = "'0.07290"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.8696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'20.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "1.871.55"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "'0.07633"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.17%  "
$ws.Range("D14").Value = "'5.335"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "'92.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "'1.011"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "'0.000008631"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "27.445.01"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "'14.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "'5.208"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "2.096.73"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").Value = "'1.880"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").Value = "'151.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "'18.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D28").Value = "'2.081"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("D29").Value = "'5.094"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("D30").Value = "'116.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").Value = "'0.08898"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'0.7363"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.05%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.447"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.32%  "

$ws.Range("D36").Value = "'1.011"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'2.479"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").Value = "'1.080"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").Value = "'0.05244"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'0.01910"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").Value = "'2.927"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").Value = "'0.5186"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").Value = "'0.1624"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").Value = "'8.275"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").Value = "'0.4824"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").Value = "'10.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").Value = "'103.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "'1.632"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").Value = "'0.06268"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.90%  "
